# Enrollment_Log.xlsx - append two new enrollment records (rows 8 and 9)
# to the "Enrollment_Log" worksheet, matching the values/order introduced
# by the commit "Adding different messages to the InputException for
# different input types".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B (Enrollment_Date), D (Age) and F (Enrollment_Arm) contain values
# that Excel would otherwise auto-convert (dates / numbers) instead of
# keeping them as plain text like the rest of the log. Temporarily format
# the new rows as Text before writing the values, then restore the default
# "General" format so no residual cell styling is left behind.
$newRows = $ws.Range("A8:G9")
$newRows.NumberFormat = "@"

# Row 8
$ws.Range("A8").Value = "k"
$ws.Range("B8").Value = "2017-10-28"
$ws.Range("C8").Value = "22:01:36.640593"
$ws.Range("D8").Value = "9"
$ws.Range("E8").Value = "m"
$ws.Range("F8").Value = "-0"
$ws.Range("G8").Value = "j"

# Row 9
$ws.Range("A9").Value = "j"
$ws.Range("B9").Value = "2017-10-28"
$ws.Range("C9").Value = "22:03:54.603169"
$ws.Range("D9").Value = "9"
$ws.Range("E9").Value = "m"
$ws.Range("F9").Value = "a"
$ws.Range("G9").Value = "d"

$newRows.NumberFormat = "General"
